$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename card_type values before shifting columns:
# CharacterCard -> Card ; TextInput -> TextInputScreen
$ws.Range("C2").Value = "Card"
$ws.Range("C3").Value = "TextInputScreen"
$ws.Range("C4").Value = "TextInputScreen"
$ws.Range("C5").Value = "TextInputScreen"

# Remove the now-unused "id" and "order" columns (A, B); everything shifts left.
$ws.Columns("A:B").Delete()

# Update the selection to reflect where editing left off.
$ws.Range("A7").Select()

# Move the saved window position (workbookView xWindow/yWindow).
$excel.Left = 18540
$excel.Top = 6060
$win = $wb.Windows.Item(1)
$win.Left = 18540
$win.Top = 6060
